# Apply crypto price/volume updates per the commit diff.
# Values are written with a leading apostrophe to force text
# interpretation (so numeric-looking strings like "212.49" or
# "1.00" keep their exact text form instead of becoming floats),
# then the style is reset to Normal so no stray NumberFormat is
# left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "26.210.53"
Set-TextValue "E2" "  -0.34%  "

Set-TextValue "D3" "1.592.70"
Set-TextValue "E3" "  +0.24%  "

Set-TextValue "E4" "  -0.08%  "

Set-TextValue "D5" "212.49"
Set-TextValue "E5" "  +1.14%  "

Set-TextValue "D6" "0.501"
Set-TextValue "E6" "  -0.52%  "

Set-TextValue "E7" "  -0.08%  "

Set-TextValue "D8" "0.246"
Set-TextValue "E8" "  -0.17%  "

Set-TextValue "D9" "0.0607"
Set-TextValue "E9" "  -0.54%  "

Set-TextValue "D10" "19.05"
Set-TextValue "E10" "  -1.56%  "

Set-TextValue "D11" "0.0849"
Set-TextValue "E11" "  +0.37%  "

Set-TextValue "E12" "  +0.28%  "

Set-TextValue "D13" "1.596.70"
Set-TextValue "E13" "  +0.44%  "

Set-TextValue "E14" "  -1.62%  "

Set-TextValue "D15" "0.510"
Set-TextValue "E15" "  -1.99%  "

Set-TextValue "D16" "63.82"
Set-TextValue "E16" "  -0.83%  "

Set-TextValue "D17" "26.216.41"
Set-TextValue "E17" "  -0.30%  "

Set-TextValue "E18" "  -0.55%  "

Set-TextValue "B19" "BitcoinCash"
Set-TextValue "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "216.13"
Set-TextValue "E19" "  +2.50%  "

Set-TextValue "B20" "Chainlink"
Set-TextValue "C20" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue "D20" "7.37"
Set-TextValue "E20" "  -0.85%  "

Set-TextValue "E21" "  -0.09%  "

Set-TextValue "D22" "4.27"
Set-TextValue "E22" "  -0.12%  "

Set-TextValue "D23" "9.06"
Set-TextValue "E23" "  +1.47%  "

Set-TextValue "D24" "2.12"
Set-TextValue "E24" "  -1.15%  "

Set-TextValue "D25" "144.98"
Set-TextValue "E25" "  +0.32%  "

Set-TextValue "E26" "  -0.18%  "

Set-TextValue "D27" "6.97"
Set-TextValue "E27" "  -0.89%  "

Set-TextValue "E28" "  -1.04%  "

Set-TextValue "D29" "15.13"
Set-TextValue "E29" "  -0.65%  "

Set-TextValue "D30" "0.0493"
Set-TextValue "E30" "  -2.23%  "

Set-TextValue "E31" "  +0.73%  "

Set-TextValue "E32" "  -1.11%  "

Set-TextValue "D33" "1.425.72"
Set-TextValue "E33" "  +8.20%  "

Set-TextValue "E34" "  -1.02%  "

Set-TextValue "E35" "  -0.96%  "

Set-TextValue "D36" "1.47"
Set-TextValue "E36" "  -0.36%  "

Set-TextValue "D37" "0.588"
Set-TextValue "E37" "  -2.63%  "

Set-TextValue "D38" "0.0166"
Set-TextValue "E38" "  -0.99%  "

Set-TextValue "D39" "0.828"
Set-TextValue "E39" "  +2.37%  "

Set-TextValue "D40" "5.91"
Set-TextValue "E40" "  +4.76%  "

Set-TextValue "D41" "1.00"
Set-TextValue "E41" "  -0.08%  "

Set-TextValue "D42" "0.978"
Set-TextValue "E42" "  -10.29%  "

Set-TextValue "D43" "0.766"
Set-TextValue "E43" "  -0.04%  "

Set-TextValue "D44" "2.14"
Set-TextValue "E44" "  -0.02%  "

Set-TextValue "D45" "1.729.18"
Set-TextValue "E45" "  +0.21%  "

Set-TextValue "D46" "61.13"
Set-TextValue "E46" "  -1.42%  "

Set-TextValue "D47" "86.98"
Set-TextValue "E47" "  -0.69%  "

Set-TextValue "E48" "  +0.50%  "

Set-TextValue "D49" "0.0501"
Set-TextValue "E49" "  -0.72%  "

Set-TextValue "D50" "0.0953"
Set-TextValue "E50" "  -2.16%  "

Set-TextValue "E51" "  -0.19%  "

